$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 24 formatting: it becomes a "bordered" style row (was plain) ---
# Copy the border/font style used by row 21 (style ids 8/9) onto row 24,
# preserving its existing values.
$ws.Range("A21:E21").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)

# --- Add new row 25 with the plain style (ids 4/5), matching what row 24 used to have ---
$ws.Range("A10:E10").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)

# Fill in the new row's data
$ws.Range("A25").Value = "SCRIPT/G01P03A/us2307.ssb"
$ws.Range("B25").Value = 19
$ws.Range("C25").Value = " Exploring is always fun! ♪"
$ws.Range("D25").Value = " Исследовать всегда весело! ♪"
$ws.Range("E25").Value = " Éòòìåäïâàóû âòåãäà âåòåìï! ♪"

# Match the same wrapped-text row height used by the other entries
$ws.Rows.Item(25).RowHeight = 43.2

# --- Update the view so the newly added row is in focus ---
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("D25").Select()
